# The workbook's first sheet has four "survey" questions (Alain / Henri /
# Tony / Dulcinee) repeated cyclically across many 4-column blocks, followed
# by an email-address column and a trailing blank column. This change adds
# eight more repeats of that 4-column block (32 new columns total) right
# before the email column, pushing the email / blank columns further right.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$blockWidth = 4
$blockCount = 8             # 8 * 4 = 32 new columns

$insertColLetter = "CW"     # first column to push right (old email column)
$lastNewColLetter = "EB"    # last of the 32 freshly inserted columns
$lastRow = 9

# 1) Insert 32 blank columns at CW:EB, shifting the existing email (CW) and
#    trailing blank (CX) columns to EC/ED.
$ws.Range("$($insertColLetter)1:$($lastNewColLetter)$lastRow").EntireColumn.Insert()

# 2) The 4 columns immediately to the left of the insertion point (CS:CV)
#    hold one full cycle of the repeating Alain/Henri/Tony/Dulcinee pattern.
#    Copy that block and paste it 8 times across the newly inserted columns
#    so the cyclic pattern simply continues.
$templateFirstCol = $ws.Range($insertColLetter + "1").Column - $blockWidth
$template = $ws.Range($ws.Cells.Item(1, $templateFirstCol), $ws.Cells.Item($lastRow, $templateFirstCol + $blockWidth - 1))
$template.Copy()

$startCol = $ws.Range($insertColLetter + "1").Column
for ($i = 0; $i -lt $blockCount; $i++) {
    $col1 = $startCol + ($i * $blockWidth)
    $col2 = $col1 + $blockWidth - 1
    $target = $ws.Range($ws.Cells.Item(1, $col1), $ws.Cells.Item($lastRow, $col2))
    $target.PasteSpecial(-4104)
}

$excel.CutCopyMode = 0
